$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row (dimension is A1:H45, but detect dynamically)
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 4).Value2  # column D: nativity / zTOTAL marker
    if ($label -eq "zTOTAL") {
        $total = $ws.Cells.Item($r, 7).Value2  # column G total for the group
        $ws.Cells.Item($r, 8).Value2 = $total
    } else {
        # Peek ahead within the current group to find the zTOTAL row's G value
        $g = $r
        while ($ws.Cells.Item($g, 4).Value2 -ne "zTOTAL" -and $g -le $lastRow) {
            $g++
        }
        $total = $ws.Cells.Item($g, 7).Value2
        $ws.Cells.Item($r, 8).Value2 = $total
    }
}
